# Added testcase for different signings for experiment.
#
# - Appends a new username/password credential row (row 3) to the "users"
#   sheet, mirroring the existing row 2 pattern: two mailto: hyperlinked
#   cells styled with the built-in "Hyperlink" style.
# - Moves the active tab/selection from "SigningExperiment" back to "users"
#   (at cell B7).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("users")
$ws2 = $wb.Worksheets.Item("SigningExperiment")

# New credential row on the "users" sheet.
$ws1.Range("A3").Value = "jestchecking@check.com"
$ws1.Range("B3").Value = "Password@123"

$ws1.Hyperlinks.Add($ws1.Range("A3"), "mailto:jestchecking@check.com")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "mailto:Password@123")

# Re-apply the Hyperlink style so the new cells match the existing ones.
$ws1.Range("A3").Style = "Hyperlink"
$ws1.Range("B3").Style = "Hyperlink"

# Active sheet/selection moves to "users"!B7 (away from "SigningExperiment").
$ws1.Activate()
$ws1.Range("B7").Select()
